$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Trim the trailing sentence from the train/test split paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This data will be split into training and test sets (70/30 split) and will be seeded for consistency of the split. The classification models will be trained on the same exact training data and will be tested with the same test dataset. Then during the results, I can do some anecdotal tests for messages I create to test out our models.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This data will be split into training and test sets (70/30 split) and will be seeded for consistency of the split. The classification models will be trained on the same exact training data and will be tested with the same test dataset.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2. Add a <w:lastRenderedPageBreak/> right before the text of the
#    "Machine Learning Model" heading run (same run, keeps its rPr/underline).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Machine Learning Model", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target1 = $d.Range($r1.Start, $r1.End)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4EA07493" w14:textId="3BCD97EA" w:rsidR="00C43D95" w:rsidRDefault="00C43D95" w:rsidP="00F575C3"><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Machine Learning Model</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target1.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# 3. Add a <w:lastRenderedPageBreak/> before "actual(ham)" in the Naive
#    Bayes confusion-matrix table (third occurrence of "actual(ham)").
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("actual(ham)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Find.Execute("actual(ham)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Find.Execute("actual(ham)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target2 = $d.Range($r2.Start, $r2.End)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="34226328" w14:textId="77777777" w:rsidR="00D27F41" w:rsidRDefault="00D27F41" w:rsidP="00FA236A"><w:r><w:lastRenderedPageBreak/><w:t>actual(ham)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target2.InsertXML($xml2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Conclusion paragraph: remove the mid-sentence <w:lastRenderedPageBreak/>
#    and merge the two runs it split back into a single run.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("In terms of actual accuracy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Expand(4) | Out-Null
$target3 = $d.Range($r3.Start, $r3.End)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="08235EBC" w14:textId="77777777" w:rsidR="007D3BF3" w:rsidRDefault="007D3BF3" w:rsidP="00F575C3"><w:r><w:t xml:space="preserve">In terms of actual accuracy, Linear SVC and Logistic Regression were the most accurate at predicting the correct classification. In retrospect, based on the class distribution of ham/spam having the &#8216;spam&#8217; class as positive makes the precision measure of this study a bit overrated. For example, at first look without seeing the confusion matrix, precision for the Random Forest model looks awesome at 100%, but to notice there were only 6 True Positives and 0 False Positives, the measure isn&#8217;t very valuable anymore. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target3.InsertXML($xml3) | Out-Null
